{"js": "// Update the date label and the 25 three-digit x one-digit multiplication\n// problems in the practice-sheet table to the next day's values.\nconst replacements = [\n  [\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"],\n  [\"895\u00d79=8055\", \"741\u00d73=2223\"],\n  [\"464\u00d77=3248\", \"299\u00d77=2093\"],\n  [\"695\u00d73=2085\", \"556\u00d74=2224\"],\n  [\"228\u00d78=1824\", \"758\u00d74=3032\"],\n  [\"363\u00d75=1815\", \"519\u00d78=4152\"],\n  [\"627\u00d79=5643\", \"633\u00d74=2532\"],\n  [\"769\u00d74=3076\", \"941\u00d76=5646\"],\n  [\"598\u00d79=5382\", \"541\u00d73=1623\"],\n  [\"840\u00d79=7560\", \"716\u00d72=1432\"],\n  [\"221\u00d74=884\", \"834\u00d75=4170\"],\n  [\"341\u00d74=1364\", \"992\u00d75=4960\"],\n  [\"713\u00d78=5704\", \"250\u00d79=2250\"],\n  [\"433\u00d73=1299\", \"539\u00d72=1078\"],\n  [\"413\u00d73=1239\", \"513\u00d78=4104\"],\n  [\"717\u00d76=4302\", \"568\u00d72=1136\"],\n  [\"412\u00d78=3296\", \"575\u00d78=4600\"],\n  [\"211\u00d73=633\", \"622\u00d79=5598\"],\n  [\"199\u00d74=796\", \"500\u00d74=2000\"],\n  [\"567\u00d79=5103\", \"594\u00d76=3564\"],\n  [\"942\u00d75=4710\", \"949\u00d78=7592\"],\n  [\"754\u00d79=6786\", \"350\u00d75=1750\"],\n  [\"544\u00d78=4352\", \"937\u00d77=6559\"],\n  [\"334\u00d74=1336\", \"966\u00d77=6762\"],\n  [\"847\u00d77=5929\", \"425\u00d75=2125\"],\n  [\"836\u00d73=2508\", \"857\u00d72=1714\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and the 25 three-digit x one-digit multiplication\n# problems in the practice-sheet table to the next day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-29 Monday\", \"2025-12-30 Tuesday\"),\n    @(\"895\u00d79=8055\", \"741\u00d73=2223\"),\n    @(\"464\u00d77=3248\", \"299\u00d77=2093\"),\n    @(\"695\u00d73=2085\", \"556\u00d74=2224\"),\n    @(\"228\u00d78=1824\", \"758\u00d74=3032\"),\n    @(\"363\u00d75=1815\", \"519\u00d78=4152\"),\n    @(\"627\u00d79=5643\", \"633\u00d74=2532\"),\n    @(\"769\u00d74=3076\", \"941\u00d76=5646\"),\n    @(\"598\u00d79=5382\", \"541\u00d73=1623\"),\n    @(\"840\u00d79=7560\", \"716\u00d72=1432\"),\n    @(\"221\u00d74=884\", \"834\u00d75=4170\"),\n    @(\"341\u00d74=1364\", \"992\u00d75=4960\"),\n    @(\"713\u00d78=5704\", \"250\u00d79=2250\"),\n    @(\"433\u00d73=1299\", \"539\u00d72=1078\"),\n    @(\"413\u00d73=1239\", \"513\u00d78=4104\"),\n    @(\"717\u00d76=4302\", \"568\u00d72=1136\"),\n    @(\"412\u00d78=3296\", \"575\u00d78=4600\"),\n    @(\"211\u00d73=633\", \"622\u00d79=5598\"),\n    @(\"199\u00d74=796\", \"500\u00d74=2000\"),\n    @(\"567\u00d79=5103\", \"594\u00d76=3564\"),\n    @(\"942\u00d75=4710\", \"949\u00d78=7592\"),\n    @(\"754\u00d79=6786\", \"350\u00d75=1750\"),\n    @(\"544\u00d78=4352\", \"937\u00d77=6559\"),\n    @(\"334\u00d74=1336\", \"966\u00d77=6762\"),\n    @(\"847\u00d77=5929\", \"425\u00d75=2125\"),\n    @(\"836\u00d73=2508\", \"857\u00d72=1714\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1        # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        Write-Output \"WARNING: replacement not found for '$oldText'\"\n    }\n}\n"}
